$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.606.28'
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").Value = '2.607.72'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Formula = '="508.73"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -4.08%  '
$ws.Range("D6").Formula = '="146.27"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -6.07%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.43%  '
$ws.Range("D9").Value = '2.628.12'
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("D10").Formula = '="6.41"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -1.34%  '
$ws.Range("D11").Formula = '="0.104"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -4.40%  '
$ws.Range("D12").Formula = '="0.337"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -4.01%  '
$ws.Range("D13").Formula = '="0.127"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").Value = '3.064.52'
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("D15").Value = '57.967.20'
$ws.Range("E15").Value = '  -5.01%  '
$ws.Range("D16").Formula = '="21.17"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -4.31%  '
$ws.Range("E17").Value = '  -3.42%  '
$ws.Range("D18").Value = '2.614.31'
$ws.Range("E18").Value = '  -1.94%  '
$ws.Range("D19").Formula = '="4.58"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -4.17%  '
$ws.Range("D20").Formula = '="345.05"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("D22").Formula = '="6.16"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D23").Formula = '="0.998"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Formula = '="60.75"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("D25").Formula = '="0.422"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -2.09%  '
$ws.Range("D26").Value = '2.714.06'
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").Formula = '="0.995"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  -4.94%  '
$ws.Range("D29").Value = '0.0₃0821'
$ws.Range("E29").Value = '  -4.27%  '
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("D31").Formula = '="0.998"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("D33").Formula = '="18.98"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("E34").Value = '  -4.50%  '
$ws.Range("D35").Formula = '="149.08"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").Formula = '="0.994"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +11.10%  '
$ws.Range("E37").Value = '  -2.75%  '
$ws.Range("E38").Value = '  -4.55%  '
$ws.Range("D39").Formula = '="0.864"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -6.06%  '
$ws.Range("D40").Formula = '="36.27"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Formula = '="1.42"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -4.74%  '
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Formula = '="290.51"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -5.00%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Formula = '="0.618"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -4.24%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Formula = '="0.0997"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Formula = '="19.64"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").Formula = '="0.0539"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -4.39%  '
$ws.Range("D49").Formula = '="4.74"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -3.51%  '
$ws.Range("E50").Value = '  -5.05%  '
$ws.Range("D51").Formula = '="10.26"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -1.08%  '
